# Update tuning data with 3 new sweep runs (DeepLabV3, DeepLabV3+, PAN).
# Rows 20-34 (old) become rows 22-24 and 26-37 (shifted down), while two
# brand-new rows are inserted at the top (20,21) and one more new row is
# inserted at position 25, per the commit's recorded sweep results.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the old A20:O34 block (15 rows x 15 cols) before overwriting it.
$old = $ws.Range("A20:O34").Value()

# The three brand-new sweep rows introduced by this commit.
$new1 = @(434, "pretty-sweep-24",   "finished", 7004, '["Lumen"]', "DeepLabV3",     "resnet50",         896, "Adam",    0.00001, 0.988914430141449, 0.9781160652637479, 0.990272849798202, 0.991412878036499, 9)
$new2 = @(388, "splendid-sweep-70", "finished", 3559, '["Lumen"]', "DeepLabV3Plus", "resnet101",        512, "RMSprop", 0.00001, 0.988765001296997, 0.97782814502716,   0.990845054388046, 0.989256829023361, 8)
$new3 = @(389, "hearty-sweep-69",   "finished", 3593, '["Lumen"]', "PAN",           "efficientnet-b0",  768, "Adam",    0.0001,  0.988418132066726, 0.977133482694626, 0.991938829421997, 0.989815771579742, 8)

# Build the final 18-row (20..37) x 15-col block.
$final = New-Object 'object[,]' 18,15

# For target row index 0..17 (sheet rows 20..37), -1 means "take from $oldRowMap"
# (1-based row within the old 15x15 block), 0 means brand-new row (see $newRowMap).
$oldRowMap = @(0, 0, 1, 2, 3, 0, 4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15)
$newRowMap = @($new1, $new2, $null, $null, $null, $new3, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null)

for ($r = 0; $r -lt 18; $r++) {
    $srcOldRow = $oldRowMap[$r]
    if ($srcOldRow -gt 0) {
        for ($c = 1; $c -le 15; $c++) {
            $final[$r, $c - 1] = $old[$srcOldRow, $c]
        }
    }
    else {
        $rowVals = $newRowMap[$r]
        for ($c = 0; $c -lt 15; $c++) {
            $final[$r, $c] = $rowVals[$c]
        }
    }
}

$ws.Range("A20:O37").Value = $final
